# Auto-generated edit script applying cached-value updates to Table_* sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(9, 8).Value = 7917.769
$ws.Cells.Item(9, 9).Value = 8548.416999999999
$ws.Cells.Item(9, 10).Value = 350
$ws.Cells.Item(9, 11).Value = 8548.416999999999
$ws.Cells.Item(9, 12).Value = 350
$ws.Cells.Item(9, 13).Value = -8379.416999999999
$ws.Cells.Item(9, 14).Value = -688
$ws.Cells.Item(19, 8).Value = 991
$ws.Cells.Item(19, 9).Value = 835.2222
$ws.Cells.Item(19, 10).Value = 1271.4
$ws.Cells.Item(19, 11).Value = 835.2222
$ws.Cells.Item(19, 12).Value = 1271.4
$ws.Cells.Item(19, 13).Value = -660.2222
$ws.Cells.Item(19, 14).Value = -1621.4
$ws.Cells.Item(28, 8).Value = 12348317
$ws.Cells.Item(28, 9).Value = 15875175
$ws.Cells.Item(28, 11).Value = 15875175
$ws.Cells.Item(28, 13).Value = -15874690
$ws.Cells.Item(40, 8).Value = 3242.25
$ws.Cells.Item(40, 9).Value = 2974.5
$ws.Cells.Item(40, 10).Value = 3331.5
$ws.Cells.Item(40, 11).Value = 2974.5
$ws.Cells.Item(40, 12).Value = 3331.5
$ws.Cells.Item(40, 13).Value = -2799.5
$ws.Cells.Item(40, 14).Value = -3681.5
$ws.Cells.Item(43, 8).Value = 1997
$ws.Cells.Item(43, 9).Value = 0
$ws.Cells.Item(43, 11).Value = 0
$ws.Cells.Item(43, 13).ClearContents()
$ws.Cells.Item(51, 8).Value = 7420.6562
$ws.Cells.Item(51, 9).Value = 8272.637000000001
$ws.Cells.Item(51, 10).Value = 6974.381
$ws.Cells.Item(51, 11).Value = 8272.637000000001
$ws.Cells.Item(51, 12).Value = 6974.381
$ws.Cells.Item(51, 13).Value = -7788.637000000001
$ws.Cells.Item(51, 14).Value = -7942.381
$ws.Cells.Item(62, 9).Value = 0
$ws.Cells.Item(62, 11).Value = 0
$ws.Cells.Item(62, 13).ClearContents()
$ws.Cells.Item(65, 9).Value = 0
$ws.Cells.Item(65, 11).Value = 0
$ws.Cells.Item(65, 13).ClearContents()
$ws.Cells.Item(70, 8).Value = 11789.182
$ws.Cells.Item(70, 9).Value = 2996.4
$ws.Cells.Item(70, 11).Value = 8989.200000000001
$ws.Cells.Item(70, 13).Value = -8719.200000000001
$ws.Cells.Item(73, 8).Value = 11789.182
$ws.Cells.Item(73, 9).Value = 2996.4
$ws.Cells.Item(73, 11).Value = 8989.200000000001
$ws.Cells.Item(73, 13).Value = -8053.200000000001
$ws.Cells.Item(86, 8).Value = 4245.3335
$ws.Cells.Item(86, 9).Value = 3867.1667
$ws.Cells.Item(86, 10).Value = 5001.6665
$ws.Cells.Item(86, 11).Value = 3867.1667
$ws.Cells.Item(86, 12).Value = 5001.6665
$ws.Cells.Item(86, 13).Value = -2744.1667
$ws.Cells.Item(86, 14).Value = -7247.6665
$ws.Cells.Item(89, 8).Value = 4245.3335
$ws.Cells.Item(89, 9).Value = 3867.1667
$ws.Cells.Item(89, 10).Value = 5001.6665
$ws.Cells.Item(89, 11).Value = 19335.8335
$ws.Cells.Item(89, 12).Value = 25008.3325
$ws.Cells.Item(89, 13).Value = -13719.8335
$ws.Cells.Item(89, 14).Value = -36240.3325
$ws.Cells.Item(111, 8).Value = 4329.6665
$ws.Cells.Item(111, 10).Value = 4999
$ws.Cells.Item(111, 12).Value = 14997
$ws.Cells.Item(111, 14).Value = -21131
$ws.Cells.Item(112, 8).Value = 4664.7417
$ws.Cells.Item(112, 10).Value = 4892.577
$ws.Cells.Item(112, 12).Value = 14677.731
$ws.Cells.Item(112, 14).Value = -16893.731
$ws.Cells.Item(116, 8).Value = 4148
$ws.Cells.Item(116, 9).Value = 4148
$ws.Cells.Item(116, 11).Value = 4148
$ws.Cells.Item(116, 13).Value = -706
$ws.Cells.Item(118, 8).Value = 1601.8889
$ws.Cells.Item(118, 9).Value = 1302.125
$ws.Cells.Item(118, 10).Value = 4000
$ws.Cells.Item(118, 11).Value = 3906.375
$ws.Cells.Item(118, 12).Value = 12000
$ws.Cells.Item(118, 13).Value = -2249.375
$ws.Cells.Item(118, 14).Value = -15314
$ws.Cells.Item(120, 8).Value = 21760
$ws.Cells.Item(120, 10).Value = 21760
$ws.Cells.Item(120, 12).Value = 21760
$ws.Cells.Item(120, 14).Value = -31436
$ws.Cells.Item(135, 8).Value = 2803.0386
$ws.Cells.Item(135, 9).Value = 2465.8
$ws.Cells.Item(135, 11).Value = 22192.2
$ws.Cells.Item(135, 13).Value = -19657.2
$ws.Cells.Item(137, 8).Value = 1382.0834
$ws.Cells.Item(137, 9).Value = 971.5714
$ws.Cells.Item(137, 11).Value = 2914.7142
$ws.Cells.Item(137, 13).Value = -364.7142000000003

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 841
$ws.Cells.Item(2, 9).Value = 680
$ws.Cells.Item(2, 11).Value = 680
$ws.Cells.Item(2, 13).Value = -567
$ws.Cells.Item(32, 8).Value = 2787.6482
$ws.Cells.Item(32, 9).Value = 2254.7114
$ws.Cells.Item(32, 11).Value = 2254.7114
$ws.Cells.Item(32, 13).Value = -1967.7114
$ws.Cells.Item(45, 8).Value = 2358.4375
$ws.Cells.Item(45, 9).Value = 2158.7
$ws.Cells.Item(45, 10).Value = 2691.3333
$ws.Cells.Item(45, 11).Value = 2158.7
$ws.Cells.Item(45, 12).Value = 2691.3333
$ws.Cells.Item(45, 13).Value = -1781.7
$ws.Cells.Item(45, 14).Value = -3445.3333
$ws.Cells.Item(46, 8).Value = 2833.3333
$ws.Cells.Item(46, 10).Value = 2625
$ws.Cells.Item(46, 12).Value = 2625
$ws.Cells.Item(46, 14).Value = -3263
$ws.Cells.Item(88, 8).Value = 7950.4
$ws.Cells.Item(88, 10).Value = 13776.75
$ws.Cells.Item(88, 12).Value = 13776.75
$ws.Cells.Item(88, 14).Value = -14588.75
$ws.Cells.Item(91, 8).Value = 7950.4
$ws.Cells.Item(91, 10).Value = 13776.75
$ws.Cells.Item(91, 12).Value = 13776.75
$ws.Cells.Item(91, 14).Value = -16584.75
$ws.Cells.Item(116, 8).Value = 841
$ws.Cells.Item(116, 9).Value = 680
$ws.Cells.Item(116, 11).Value = 680
$ws.Cells.Item(116, 13).Value = 1614
$ws.Cells.Item(135, 8).Value = 82899.39999999999
$ws.Cells.Item(135, 10).Value = 82899.39999999999
$ws.Cells.Item(135, 12).Value = 82899.39999999999
$ws.Cells.Item(135, 14).Value = -93039.39999999999
$ws.Cells.Item(139, 8).Value = 61536.8
$ws.Cells.Item(139, 10).Value = 61536.8
$ws.Cells.Item(139, 12).Value = 61536.8
$ws.Cells.Item(139, 14).Value = -71816.8

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 841
$ws.Cells.Item(3, 9).Value = 680
$ws.Cells.Item(3, 11).Value = 680
$ws.Cells.Item(3, 13).Value = -566
$ws.Cells.Item(20, 8).Value = 2119.3333
$ws.Cells.Item(20, 9).Value = 1951.75
$ws.Cells.Item(20, 10).Value = 2454.5
$ws.Cells.Item(20, 11).Value = 1951.75
$ws.Cells.Item(20, 12).Value = 2454.5
$ws.Cells.Item(20, 13).Value = -1704.75
$ws.Cells.Item(20, 14).Value = -2948.5
$ws.Cells.Item(99, 8).Value = 1549.1333
$ws.Cells.Item(99, 9).Value = 689.4167
$ws.Cells.Item(99, 11).Value = 689.4167
$ws.Cells.Item(99, 13).Value = 808.5833

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(7, 8).Value = 124.388885
$ws.Cells.Item(7, 9).Value = 173.57143
$ws.Cells.Item(7, 10).Value = 93.09090999999999
$ws.Cells.Item(7, 11).Value = 173.57143
$ws.Cells.Item(7, 12).Value = 93.09090999999999
$ws.Cells.Item(7, 13).Value = -60.57142999999999
$ws.Cells.Item(7, 14).Value = -319.09091
$ws.Cells.Item(99, 8).Value = 1856.7142
$ws.Cells.Item(99, 9).Value = 1749.5
$ws.Cells.Item(99, 11).Value = 1749.5
$ws.Cells.Item(99, 13).Value = -251.5
$ws.Cells.Item(107, 8).Value = 649.76
$ws.Cells.Item(107, 9).Value = 548.58826
$ws.Cells.Item(107, 10).Value = 864.75
$ws.Cells.Item(107, 11).Value = 548.58826
$ws.Cells.Item(107, 12).Value = 864.75
$ws.Cells.Item(107, 13).Value = 1371.41174
$ws.Cells.Item(107, 14).Value = -4704.75
$ws.Cells.Item(126, 8).Value = 1856.7142
$ws.Cells.Item(126, 9).Value = 1749.5
$ws.Cells.Item(126, 11).Value = 5248.5
$ws.Cells.Item(126, 13).Value = -2778.5
$ws.Cells.Item(135, 8).Value = 89998
$ws.Cells.Item(135, 10).Value = 89998
$ws.Cells.Item(135, 12).Value = 89998
$ws.Cells.Item(135, 14).Value = -100138
$ws.Cells.Item(140, 8).Value = 139197.5
$ws.Cells.Item(140, 10).Value = 139197.5
$ws.Cells.Item(140, 12).Value = 139197.5
$ws.Cells.Item(140, 14).Value = -149557.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(33, 8).Value = 1679.7142
$ws.Cells.Item(33, 10).Value = 436.33334
$ws.Cells.Item(33, 12).Value = 2618.00004
$ws.Cells.Item(33, 14).Value = -3184.00004

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(10, 8).Value = 186666
$ws.Cells.Item(10, 10).Value = 29999
$ws.Cells.Item(10, 12).Value = 29999
$ws.Cells.Item(10, 14).Value = -30337
$ws.Cells.Item(122, 8).Value = 2989.3489
$ws.Cells.Item(122, 9).Value = 2806.423
$ws.Cells.Item(122, 10).Value = 3269.1177
$ws.Cells.Item(122, 11).Value = 8419.269
$ws.Cells.Item(122, 12).Value = 9807.3531
$ws.Cells.Item(122, 13).Value = -5969.269
$ws.Cells.Item(122, 14).Value = -14707.3531

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 1435.6552
$ws.Cells.Item(22, 9).Value = 1262.4286
$ws.Cells.Item(22, 11).Value = 1262.4286
$ws.Cells.Item(22, 13).Value = -967.4286
$ws.Cells.Item(27, 8).Value = 1435.6552
$ws.Cells.Item(27, 9).Value = 1262.4286
$ws.Cells.Item(27, 11).Value = 1262.4286
$ws.Cells.Item(27, 13).Value = -1155.4286
$ws.Cells.Item(40, 8).Value = 2545.5715
$ws.Cells.Item(40, 9).Value = 1645
$ws.Cells.Item(40, 11).Value = 1645
$ws.Cells.Item(40, 13).Value = -1509
$ws.Cells.Item(46, 8).Value = 1847.0834
$ws.Cells.Item(46, 9).Value = 1200.3334
$ws.Cells.Item(46, 10).Value = 2062.6667
$ws.Cells.Item(46, 11).Value = 1200.3334
$ws.Cells.Item(46, 12).Value = 2062.6667
$ws.Cells.Item(46, 13).Value = -1012.3334
$ws.Cells.Item(46, 14).Value = -2438.6667
$ws.Cells.Item(93, 8).Value = 1032.4762
$ws.Cells.Item(93, 9).Value = 901.1539
$ws.Cells.Item(93, 11).Value = 901.1539
$ws.Cells.Item(93, 13).Value = 346.8461

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(122, 8).Value = 6717.476
$ws.Cells.Item(122, 9).Value = 2128.3684
$ws.Cells.Item(122, 11).Value = 6385.1052
$ws.Cells.Item(122, 13).Value = -3935.1052
$ws.Cells.Item(141, 8).Value = 94892.11
$ws.Cells.Item(141, 10).Value = 94892.11
$ws.Cells.Item(141, 12).Value = 94892.11
$ws.Cells.Item(141, 14).Value = -105252.11

